$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell V2 value from True to False
$ws.Range("V2").Value = "False"

# Update the active window's scroll position and selection
$excel.ActiveWindow.ScrollColumn = 11
$ws.Range("V2").Select()
